$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.095077037811279
$ws.Range("B1").Value = 2.726522207260132
$ws.Range("C1").Value = 2.864420413970947
$ws.Range("D1").Value = 2.862509250640869
$ws.Range("E1").Value = 0.7830931544303894
